$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day 3 new male pair: fill in D column for rows 7 and 8 (coop014 / coop015),
# and their Average (F) / Weight (G) formulas.
$ws.Range("D7").Value = 28.2
$ws.Range("F7").Formula = "=AVERAGE(B7:D7)"
$ws.Range("G7").Formula = "=F7*0.8"
$ws.Range("F7:G7").NumberFormat = "0.0"

$ws.Range("D8").Value = 29.4
$ws.Range("F8").Formula = "=AVERAGE(B8:D8)"
# Written together with the existing G3:G5 weight-formula range so the
# G column keeps sharing that formula pattern down through G8.
$ws.Range("G8,G3:G5").Formula = "=F8*0.8"
$ws.Range("F8:G8").NumberFormat = "0.0"

# Update selection to G8
$ws.Range("G8").Select()
